$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1754.3334
$ws.Range("J62").Value = 2044.75
$ws.Range("L62").Value = 2044.75
$ws.Range("N62").Value = -3292.75
$ws.Range("H65").Value = 1754.3334
$ws.Range("J65").Value = 2044.75
$ws.Range("L65").Value = 10223.75
$ws.Range("N65").Value = -16463.75
$ws.Range("H132").Value = 2467.5508
$ws.Range("I132").Value = 2354.9524
$ws.Range("K132").Value = 7064.8572
$ws.Range("M132").Value = -4534.8572
$ws.Range("H138").Value = 2032.2812
$ws.Range("I138").Value = 1283.3695
$ws.Range("J138").Value = 2721.28
$ws.Range("K138").Value = 3850.1085
$ws.Range("L138").Value = 8163.84
$ws.Range("M138").Value = 1289.8915
$ws.Range("N138").Value = -18443.84
$ws.Range("H139").Value = 47400
$ws.Range("J139").Value = 47400
$ws.Range("L139").Value = 47400
$ws.Range("N139").Value = -57680

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2446.3022
$ws.Range("I61").Value = 1704.3334
$ws.Range("J61").Value = 4158.5386
$ws.Range("K61").Value = 1704.3334
$ws.Range("L61").Value = 4158.5386
$ws.Range("M61").Value = -1492.3334
$ws.Range("N61").Value = -4582.5386
$ws.Range("H74").Value = 1302.8125
$ws.Range("I74").Value = 1056.3334
$ws.Range("K74").Value = 1056.3334
$ws.Range("M74").Value = -182.3334
$ws.Range("H77").Value = 1302.8125
$ws.Range("I77").Value = 1056.3334
$ws.Range("K77").Value = 5281.666999999999
$ws.Range("M77").Value = -913.6669999999995
$ws.Range("H101").Value = 78401.336
$ws.Range("J101").Value = 78401.336
$ws.Range("L101").Value = 78401.336
$ws.Range("N101").Value = -84891.336
$ws.Range("H112").Value = 45994.09
$ws.Range("J112").Value = 45994.09
$ws.Range("L112").Value = 45994.09
$ws.Range("N112").Value = -48948.09
$ws.Range("H122").Value = 39358.81
$ws.Range("I122").Value = 53396.79
$ws.Range("J122").Value = 1255.7142
$ws.Range("K122").Value = 160190.37
$ws.Range("L122").Value = 3767.1426
$ws.Range("M122").Value = -157740.37
$ws.Range("N122").Value = -8667.142599999999
$ws.Range("H129").Value = 52780
$ws.Range("J129").Value = 52780
$ws.Range("L129").Value = 52780
$ws.Range("N129").Value = -62780
$ws.Range("H132").Value = 2275.1792
$ws.Range("I132").Value = 1662.125
$ws.Range("J132").Value = 3823.9473
$ws.Range("K132").Value = 4986.375
$ws.Range("L132").Value = 11471.8419
$ws.Range("M132").Value = -2456.375
$ws.Range("N132").Value = -16531.8419
$ws.Range("H136").Value = 2446.3022
$ws.Range("I136").Value = 1704.3334
$ws.Range("J136").Value = 4158.5386
$ws.Range("K136").Value = 5113.0002
$ws.Range("L136").Value = 12475.6158
$ws.Range("M136").Value = -2563.0002
$ws.Range("N136").Value = -17575.6158

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 21996.4
$ws.Range("J35").Value = 21996.4
$ws.Range("L35").Value = 21996.4
$ws.Range("N35").Value = -22616.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5082.839
$ws.Range("I134").Value = 4902.3335
$ws.Range("J134").Value = 6301.25
$ws.Range("K134").Value = 14707.0005
$ws.Range("L134").Value = 18903.75
$ws.Range("M134").Value = -12172.0005
$ws.Range("N134").Value = -23973.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 11114245
$ws.Range("J32").Value = 11114245
$ws.Range("L32").Value = 33342735
$ws.Range("N32").Value = -33343301
$ws.Range("H46").Value = 1173.8889
$ws.Range("I46").Value = 456.42856
$ws.Range("J46").Value = 1425
$ws.Range("K46").Value = 1369.28568
$ws.Range("L46").Value = 4275
$ws.Range("M46").Value = -1278.28568
$ws.Range("N46").Value = -4457
$ws.Range("H131").Value = 8231.8125
$ws.Range("I131").Value = 416.27274
$ws.Range("J131").Value = 25426
$ws.Range("K131").Value = 1248.81822
$ws.Range("L131").Value = 76278
$ws.Range("M131").Value = 3791.18178
$ws.Range("N131").Value = -86358
$ws.Range("H134").Value = 4922.706
$ws.Range("I134").Value = 1973.625
$ws.Range("J134").Value = 7544.1113
$ws.Range("K134").Value = 5920.875
$ws.Range("L134").Value = 22632.3339
$ws.Range("M134").Value = -850.875
$ws.Range("N134").Value = -32772.3339
$ws.Range("H137").Value = 6179994
$ws.Range("J137").Value = 4140.5625
$ws.Range("L137").Value = 12421.6875
$ws.Range("N137").Value = -22621.6875
$ws.Range("H139").Value = 3306.195
$ws.Range("J139").Value = 4399.75
$ws.Range("L139").Value = 13199.25
$ws.Range("N139").Value = -23479.25
$ws.Range("H140").Value = 1815.8334
$ws.Range("I140").Value = 1346.5217
$ws.Range("J140").Value = 2646.1538
$ws.Range("K140").Value = 4039.5651
$ws.Range("L140").Value = 7938.4614
$ws.Range("M140").Value = 1140.4349
$ws.Range("N140").Value = -18298.4614

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 85009
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 85009
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 85009
$ws.Range("N25").Value = -86067
$ws.Range("H27").Value = 7000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 7000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 7000
$ws.Range("N27").Value = -7332
$ws.Range("H70").Value = 8848.781999999999
$ws.Range("I70").Value = 10520.25
$ws.Range("J70").Value = 5028.2856
$ws.Range("K70").Value = 10520.25
$ws.Range("L70").Value = 5028.2856
$ws.Range("M70").Value = -10250.25
$ws.Range("N70").Value = -5568.2856
$ws.Range("H73").Value = 8848.781999999999
$ws.Range("I73").Value = 10520.25
$ws.Range("J73").Value = 5028.2856
$ws.Range("K73").Value = 10520.25
$ws.Range("L73").Value = 5028.2856
$ws.Range("M73").Value = -9584.25
$ws.Range("N73").Value = -6900.2856
$ws.Range("H102").Value = 1869.6451
$ws.Range("I102").Value = 1727.4584
$ws.Range("K102").Value = 1727.4584
$ws.Range("M102").Value = -105.4584
$ws.Range("H122").Value = 1985.6
$ws.Range("I122").Value = 1415.1666
$ws.Range("J122").Value = 2841.25
$ws.Range("K122").Value = 4245.4998
$ws.Range("L122").Value = 8523.75
$ws.Range("M122").Value = -1795.4998
$ws.Range("N122").Value = -13423.75
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("M27").ClearContents()
$ws.Range("N131").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9744
$ws.Range("J54").Value = 9744
$ws.Range("L54").Value = 9744
$ws.Range("N54").Value = -10784
$ws.Range("H81").Value = 4259.231
$ws.Range("J81").Value = 3850
$ws.Range("L81").Value = 7700
$ws.Range("N81").Value = -9822
$ws.Range("H84").Value = 4259.231
$ws.Range("J84").Value = 3850
$ws.Range("L84").Value = 38500
$ws.Range("N84").Value = -49108
$ws.Range("H132").Value = 4505758
$ws.Range("I132").Value = 1242.7727
$ws.Range("J132").Value = 11112380
$ws.Range("K132").Value = 3728.3181
$ws.Range("L132").Value = 33337140
$ws.Range("M132").Value = -1198.3181
$ws.Range("N132").Value = -33342200
$ws.Range("H136").Value = 2290.2754
$ws.Range("J136").Value = 5266.5
$ws.Range("L136").Value = 15799.5
$ws.Range("N136").Value = -20899.5

Write-Output "Applied all updates"